$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.933575034141541
$ws.Range("B1").Value = 2.66878342628479
$ws.Range("C1").Value = 2.865574836730957
$ws.Range("D1").Value = 3.434018611907959
$ws.Range("E1").Value = 2.209897518157959
